$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.899.04'
$ws.Range("E2").Value = '  -1.53%  '

$ws.Range("D3").Value = '3.062.61'
$ws.Range("E3").Value = '  -1.18%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '558.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.51%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.69'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.17%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("D8").Value = '3.060.17'
$ws.Range("E8").Value = '  -1.17%  '

$ws.Range("E9").Value = '  +3.23%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.153'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.30%  '

$ws.Range("E11").Value = '  -2.93%  '

$ws.Range("E12").Value = '  +1.48%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000232'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.86%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.28'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.46%  '

$ws.Range("D15").Value = '3.569.05'
$ws.Range("E15").Value = '  -0.92%  '

$ws.Range("D16").Value = '63.965.48'
$ws.Range("E16").Value = '  -1.54%  '

$ws.Range("D17").Value = '3.061.82'
$ws.Range("E17").Value = '  -1.22%  '

$ws.Range("E18").Value = '  +0.13%  '

$ws.Range("E19").Value = '  -0.07%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '486.83'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.37%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.34'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.68%  '

$ws.Range("E22").Value = '  -0.09%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.61'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.89%  '

$ws.Range("E24").Value = '  -0.54%  '

$ws.Range("E25").Value = '  +1.89%  '

$ws.Range("E26").Value = '  +0.04%  '

$ws.Range("E27").Value = '  +0.16%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.15'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.63%  '

$ws.Range("E29").Value = '  -0.91%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.07%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.55'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.40%  '

$ws.Range("E32").Value = '  +0.25%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.51'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.11%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.69'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.89%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.23'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.03%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '54.95'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0412'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.73%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '444.28'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.99%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0815'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.80%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.80'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.36%  '

$ws.Range("D41").Value = '3.024.89'
$ws.Range("E41").Value = '  +1.28%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.33'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.79%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.277'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.28%  '

$ws.Range("E45").Value = '  +4.80%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '27.75'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.25%  '

$ws.Range("E47").Value = '  +0.00%  '

$ws.Range("E48").Value = '  +0.85%  '

$ws.Range("B49").Value = 'PEPE'
$ws.Range("C49").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D49").Value = '0.0₃0517'
$ws.Range("E49").Value = '  -3.19%  '

$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '117.79'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.16%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.13'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.39%  '
